$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "동탄숲속마을자연앤경남아너스빌(1124)" (old row 35).
# This shifts all subsequent rows up by one and Excel automatically drops
# the now-unreferenced shared string from the table.
$ws.Rows(35).Delete()

# Leave the selection where the author last clicked before saving.
$ws.Range("C20").Select()
